$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts N->O, O->P, P->Q).
# Excel copies the formatting of the column to the left (M) onto the new
# column, which is exactly what happened in the target workbook.
$ws.Columns("N").Insert()

# The new column keeps the same numeric width as column M (11 chars) but
# without the "best fit" auto-size flag, matching the target XML.
$ws.Columns("N").ColumnWidth = 10.14

# Make "Repayment schedule" the active sheet/tab and select Q6 on it
# (this also clears the tabSelected flag that used to be on "Prepay Loan").
$ws.Activate()
[void]$ws.Range("Q6").Select()
